$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Drop the unused trailing columns (E:I) that carried the old
#    "width 30" custom formatting - only A:D keep a custom width now.
#    This also clears whatever the old E1:I1 header cells held.
# ------------------------------------------------------------------
$ws.Columns("E:I").Delete()

# ------------------------------------------------------------------
# 2. Re-write the header row (row 1). Columns E1:G1 were wiped by the
#    column delete above, so re-create them and restore the bold /
#    bordered / centered header look by copying the format from A1
#    (which already carries that style) - this keeps the style table
#    from growing with near-duplicate entries.
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Fecha"
$ws.Range("B1").Value = "Codigo"
$ws.Range("C1").Value = "Cliente"
$ws.Range("D1").Value = "Nro operacion"
$ws.Range("E1").Value = "Nro operacion 2"
$ws.Range("F1").Value = "Moneda"
$ws.Range("G1").Value = "Monto sin IGV"

$ws.Range("A1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3. Columns D and G hold numbers-that-look-like-numbers but must stay
#    plain text (invoice/operation numbers, amounts kept as text).
#    Pre-mark the ranges as Text so the values below aren't silently
#    coerced into numeric cells, then strip the temporary formatting
#    back off once the text is safely stored.
# ------------------------------------------------------------------
$ws.Range("D2:D14").NumberFormat = "@"
$ws.Range("G2:G16").NumberFormat = "@"

# ------------------------------------------------------------------
# 4. Data rows.
# ------------------------------------------------------------------
$ws.Range("A2").Value = 44928
$ws.Range("B2").Value = "F001-0298"
$ws.Range("C2").Value = "PESQUERA CONCEPCION S.A.C."
$ws.Range("D2").Value = "-"
$ws.Range("F2").Value = "DOLARES"
$ws.Range("G2").Value = "2867.4"

$ws.Range("A3").Value = 44928
$ws.Range("B3").Value = "F001-0299"
$ws.Range("C3").Value = "MEGUI INVESTMENT S.A.C."
$ws.Range("D3").Value = "171693"
$ws.Range("F3").Value = "DOLARES"
$ws.Range("G3").Value = "412.14"

$ws.Range("A4").Value = 44932
$ws.Range("B4").Value = "F001-0300"
$ws.Range("C4").Value = "EMPRESA GRUPO JIVO SOCIEDAD ANONIMA CERRADA"
$ws.Range("D4").Value = "1273480"
$ws.Range("F4").Value = "DOLARES"
$ws.Range("G4").Value = "1658.88"

$ws.Range("A5").Value = 44932
$ws.Range("B5").Value = "F001-0304"
$ws.Range("C5").Value = "OROSCO CASTRO JIMMY NICOLAY"
$ws.Range("D5").Value = "69082"
$ws.Range("F5").Value = "SOLES"
$ws.Range("G5").Value = "390.3"

$ws.Range("A6").Value = 44932
$ws.Range("B6").Value = "F001-0305"
$ws.Range("C6").Value = "EMPRESA PESQUERA ROSA ISABEL S.R.L."
$ws.Range("D6").Value = "1830879"
$ws.Range("F6").Value = "DOLARES"
$ws.Range("G6").Value = "2191.7"

$ws.Range("A7").Value = 44932
$ws.Range("B7").Value = "F001-0306"
$ws.Range("C7").Value = "OROSCO CASTRO JIMMY NICOLAY"
$ws.Range("D7").Value = "1986175"
$ws.Range("F7").Value = "SOLES"
$ws.Range("G7").Value = "68.78"

$ws.Range("A8").Value = 44932
$ws.Range("B8").Value = "F001-0307"
$ws.Range("C8").Value = "STEEL ASESORIA E.I.R.L."
$ws.Range("D8").Value = "243534597"
$ws.Range("F8").Value = "SOLES"
$ws.Range("G8").Value = "146.16"

$ws.Range("A9").Value = 44932
$ws.Range("B9").Value = "F001-0308"
$ws.Range("C9").Value = "STEEL ASESORIA E.I.R.L."
$ws.Range("D9").Value = "243536416"
$ws.Range("F9").Value = "SOLES"
$ws.Range("G9").Value = "213.21"

$ws.Range("A10").Value = 44932
$ws.Range("B10").Value = "F001-0309"
$ws.Range("C10").Value = "MEGUI INVESTMENT S.A.C."
$ws.Range("D10").Value = "39411"
$ws.Range("F10").Value = "DOLARES"
$ws.Range("G10").Value = "4149.78"

$ws.Range("A11").Value = 44932
$ws.Range("B11").Value = "F001-0311"
$ws.Range("C11").Value = "EMPRESA PESQUERA ROSA ISABEL S.R.L."
$ws.Range("D11").Value = "1146102"
$ws.Range("F11").Value = "DOLARES"
$ws.Range("G11").Value = "429.74"

$ws.Range("A12").Value = 44932
$ws.Range("B12").Value = "F001-0312"
$ws.Range("C12").Value = "AITANA&KHALEESI E.I.R.L."
$ws.Range("D12").Value = "185"
$ws.Range("F12").Value = "DOLARES"
$ws.Range("G12").Value = "3186.72"

$ws.Range("A13").Value = 44932
$ws.Range("B13").Value = "F001-0313"
$ws.Range("C13").Value = "SERVICIOS INDUSTRIALES HALAVISI S.A.C."
$ws.Range("D13").Value = "1785131"
$ws.Range("F13").Value = "SOLES"
$ws.Range("G13").Value = "178.29"

$ws.Range("A14").Value = 44932
$ws.Range("B14").Value = "F001-0314"
$ws.Range("C14").Value = "SERVICIOS INDUSTRIALES HALAVISI S.A.C."
$ws.Range("D14").Value = "1756488"
$ws.Range("F14").Value = "SOLES"
$ws.Range("G14").Value = "73.16"

$ws.Range("F15").Value = "Monto total"
$ws.Range("G15").Value = "15966.26"

$ws.Range("F16").Value = "Monto comision"
$ws.Range("G16").Value = "159.66"

# ------------------------------------------------------------------
# 5. Drop the temporary Text formatting now that the values are
#    safely stored as text, so the cells fall back to the default
#    (unstyled) look, matching every other data cell in the sheet.
# ------------------------------------------------------------------
$ws.Range("D2:D14").ClearFormats()
$ws.Range("G2:G16").ClearFormats()

# ------------------------------------------------------------------
# 6. Column A holds real dates - format it accordingly.
# ------------------------------------------------------------------
$ws.Range("A2:A14").NumberFormat = "YYYY-MM-DD"
